# Fix bug in convertJsonToExcel: update stale "last_edited_time" string
# and the figures for row 7 (Tháng 8) that were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2024-08-23T13:32:00.000Z"
$newTimestamp = "2024-08-24T20:33:00.000Z"

# The "last_edited_time" column (D) shares the same string across several
# rows; replace every occurrence so the shared string table stays de-duplicated.
$used = $ws.UsedRange
foreach ($cell in $used.Columns.Item(4).Cells) {
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# Updated numeric figures on row 7 (Tháng 8)
$ws.Range("S7").Value = 69005000
$ws.Range("W7").Value = 105195000
$ws.Range("AE7").Value = 174200000
$ws.Range("AH7").Value = 157200000
$ws.Range("AK7").Value = 21
$ws.Range("AQ7").Value = 189200000
